$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 16.13
$ws.Range("F2").Value = 26.21
$ws.Range("K2").Value = 57.8
$ws.Range("N2").Value = 85.96878041621773

# Row 3
$ws.Range("K3").Value = 53
$ws.Range("N3").Value = 85.96878041621773

# Row 4
$ws.Range("D4").Value = 90734.73
$ws.Range("K4").Value = 50.8
$ws.Range("N4").Value = 85.96878041621773

# Row 5
$ws.Range("D5").Value = 11.81
$ws.Range("F5").Value = 15.33
$ws.Range("K5").Value = 50.2
$ws.Range("N5").Value = 85.96878041621773

# Row 6
$ws.Range("D6").Value = 177.18
$ws.Range("F6").Value = 0.03
$ws.Range("K6").Value = 47.8
$ws.Range("N6").Value = 85.96878041621773
